$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "296.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.18%"
$ws.Range("E2").ClearFormats()
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "11"
$ws.Range("G2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.27"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.44%"
$ws.Range("E3").ClearFormats()
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "11"
$ws.Range("G3").ClearFormats()

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.126"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.58%"
$ws.Range("E4").ClearFormats()
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "11"
$ws.Range("G4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07337"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.69%"
$ws.Range("E5").ClearFormats()
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "11"
$ws.Range("G5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.714"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.89%"
$ws.Range("E6").ClearFormats()
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "11"
$ws.Range("G6").ClearFormats()

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.735"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.00%"
$ws.Range("E7").ClearFormats()
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "11"
$ws.Range("G7").ClearFormats()

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.632"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "9.79%"
$ws.Range("E8").ClearFormats()
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "11"
$ws.Range("G8").ClearFormats()

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9190"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.56%"
$ws.Range("E9").ClearFormats()
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "11"
$ws.Range("G9").ClearFormats()

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1673"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.99%"
$ws.Range("E10").ClearFormats()
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "11"
$ws.Range("G10").ClearFormats()

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07147"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.48%"
$ws.Range("E11").ClearFormats()
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "11"
$ws.Range("G11").ClearFormats()

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07965"
$ws.Range("D12").ClearFormats()
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "11"
$ws.Range("G12").ClearFormats()

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02986"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.47%"
$ws.Range("E13").ClearFormats()
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "11"
$ws.Range("G13").ClearFormats()

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09921"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.05%"
$ws.Range("E14").ClearFormats()
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "11"
$ws.Range("G14").ClearFormats()

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001490"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.08%"
$ws.Range("E15").ClearFormats()
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "11"
$ws.Range("G15").ClearFormats()

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006100"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.36%"
$ws.Range("E16").ClearFormats()
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "11"
$ws.Range("G16").ClearFormats()

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.67%"
$ws.Range("E17").ClearFormats()
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "11"
$ws.Range("G17").ClearFormats()

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.228"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.02%"
$ws.Range("E18").ClearFormats()
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "11"
$ws.Range("G18").ClearFormats()

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3272"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.89%"
$ws.Range("E19").ClearFormats()
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "11"
$ws.Range("G19").ClearFormats()

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.42%"
$ws.Range("E20").ClearFormats()
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "11"
$ws.Range("G20").ClearFormats()

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.555"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.75%"
$ws.Range("E21").ClearFormats()
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "11"
$ws.Range("G21").ClearFormats()

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04617"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.81%"
$ws.Range("E22").ClearFormats()
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "11"
$ws.Range("G22").ClearFormats()

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1548"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.57%"
$ws.Range("E23").ClearFormats()
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "11"
$ws.Range("G23").ClearFormats()

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.31%"
$ws.Range("E24").ClearFormats()
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "11"
$ws.Range("G24").ClearFormats()

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004426"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.27%"
$ws.Range("E25").ClearFormats()
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "11"
$ws.Range("G25").ClearFormats()

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.09%"
$ws.Range("E26").ClearFormats()
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "11"
$ws.Range("G26").ClearFormats()

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "7.44%"
$ws.Range("E27").ClearFormats()
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "11"
$ws.Range("G27").ClearFormats()

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "11"
$ws.Range("G28").ClearFormats()

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "11"
$ws.Range("G29").ClearFormats()

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "11"
$ws.Range("G30").ClearFormats()

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "11"
$ws.Range("G31").ClearFormats()

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "11"
$ws.Range("G32").ClearFormats()

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "11"
$ws.Range("G33").ClearFormats()

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "11"
$ws.Range("G34").ClearFormats()

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "11"
$ws.Range("G35").ClearFormats()

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "11"
$ws.Range("G36").ClearFormats()

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "11"
$ws.Range("G37").ClearFormats()

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "11"
$ws.Range("G38").ClearFormats()

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01686"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.35%"
$ws.Range("E39").ClearFormats()
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "11"
$ws.Range("G39").ClearFormats()

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04412"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.86%"
$ws.Range("E40").ClearFormats()
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "11"
$ws.Range("G40").ClearFormats()

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007167"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.25%"
$ws.Range("E41").ClearFormats()
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "11"
$ws.Range("G41").ClearFormats()

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1329"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.99%"
$ws.Range("E42").ClearFormats()
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "11"
$ws.Range("G42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002138"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.23%"
$ws.Range("E43").ClearFormats()
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "11"
$ws.Range("G43").ClearFormats()

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01068"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-19.72%"
$ws.Range("E44").ClearFormats()
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "11"
$ws.Range("G44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006008"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.26%"
$ws.Range("E45").ClearFormats()
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "11"
$ws.Range("G45").ClearFormats()

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.88%"
$ws.Range("E46").ClearFormats()
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "11"
$ws.Range("G46").ClearFormats()

# Row 47
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "11"
$ws.Range("G47").ClearFormats()

# Row 48
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "11"
$ws.Range("G48").ClearFormats()

# Row 49
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "11"
$ws.Range("G49").ClearFormats()

# Row 50
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "11"
$ws.Range("G50").ClearFormats()

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "11"
$ws.Range("G51").ClearFormats()
